# "added 4wk low sales check"
# Updates the forecast numbers on the "Forecast Comparison" sheet (MyForecast,
# Inventory Coverage, Seasonality Index) and the corresponding roll-up
# figures on the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H), Seasonality Index (L) ---

$ws1.Range("D2").Value = 90
$ws1.Range("H2").Value = 6.72
$ws1.Range("L2").Value = 0.85

$ws1.Range("D3").Value = 90
$ws1.Range("H3").Value = 5.72

$ws1.Range("D4").Value = 89
$ws1.Range("H4").Value = 4.78
$ws1.Range("L4").Value = 1.06

$ws1.Range("D5").Value = 90
$ws1.Range("H5").Value = 3.73
$ws1.Range("L5").Value = 1.13

$ws1.Range("D6").Value = 89
$ws1.Range("H6").Value = 2.76
$ws1.Range("L6").Value = 0.96

$ws1.Range("D7").Value = 89
$ws1.Range("H7").Value = 1.76
$ws1.Range("L7").Value = 0.85

$ws1.Range("D8").Value = 90
$ws1.Range("H8").Value = 0.76
$ws1.Range("L8").Value = 0.89

$ws1.Range("D9").Value = 89
$ws1.Range("L9").Value = 1.05

$ws1.Range("D10").Value = 88
$ws1.Range("L10").Value = 1.18

$ws1.Range("D11").Value = 89
$ws1.Range("L11").Value = 0.92

$ws1.Range("D12").Value = 89
$ws1.Range("L12").Value = 1.16

$ws1.Range("D13").Value = 89
$ws1.Range("L13").Value = 1.04

$ws1.Range("D14").Value = 89
$ws1.Range("L14").Value = 1.16

$ws1.Range("D15").Value = 88
$ws1.Range("L15").Value = 1.04

$ws1.Range("D16").Value = 89
$ws1.Range("L16").Value = 0.96

$ws1.Range("D17").Value = 88
$ws1.Range("L17").Value = 1.14

# --- Summary: roll-up figures (stored as text, like the rest of column B) ---

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "1425"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "716"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "359"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "90"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "88"
